$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 936.36664
$ws.Range("I15").Value = 936.36664
$ws.Range("K15").Value = 2809.09992
$ws.Range("M15").Value = -2640.09992
$ws.Range("H69").Value = 16507.5
$ws.Range("J69").Value = 30015
$ws.Range("L69").Value = 90045
$ws.Range("N69").Value = -91793
$ws.Range("H72").Value = 16507.5
$ws.Range("J72").Value = 30015
$ws.Range("L72").Value = 270135
$ws.Range("N72").Value = -278871
$ws.Range("H107").Value = 1752.8823
$ws.Range("I107").Value = 1887.4166
$ws.Range("J107").Value = 1430
$ws.Range("K107").Value = 1887.4166
$ws.Range("L107").Value = 1430
$ws.Range("M107").Value = 32.58339999999998
$ws.Range("N107").Value = -5270
$ws.Range("H113").Value = 5365.6313
$ws.Range("I113").Value = 2299.6667
$ws.Range("J113").Value = 5940.5
$ws.Range("K113").Value = 2299.6667
$ws.Range("L113").Value = 5940.5
$ws.Range("M113").Value = 954.3332999999998
$ws.Range("N113").Value = -12448.5
$ws.Range("H132").Value = 12752640
$ws.Range("I132").Value = 16836018
$ws.Range("J132").Value = 502504.94
$ws.Range("K132").Value = 50508054
$ws.Range("L132").Value = 1507514.82
$ws.Range("M132").Value = -50505524
$ws.Range("N132").Value = -1512574.82
$ws.Range("H135").Value = 763.5
$ws.Range("I135").Value = 304.5625
$ws.Range("J135").Value = 1681.375
$ws.Range("K135").Value = 2741.0625
$ws.Range("L135").Value = 15132.375
$ws.Range("M135").Value = -206.0625
$ws.Range("N135").Value = -20202.375
$ws.Range("H138").Value = 3089.81
$ws.Range("I138").Value = 856.9706
$ws.Range("J138").Value = 4240.0605
$ws.Range("K138").Value = 2570.9118
$ws.Range("L138").Value = 12720.1815
$ws.Range("M138").Value = 2569.0882
$ws.Range("N138").Value = -23000.1815

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1461.4445
$ws.Range("I45").Value = 1353.8334
$ws.Range("J45").Value = 1676.6666
$ws.Range("K45").Value = 1353.8334
$ws.Range("L45").Value = 1676.6666
$ws.Range("M45").Value = -976.8334
$ws.Range("N45").Value = -2430.6666
$ws.Range("H61").Value = 1069.0333
$ws.Range("I61").Value = 915.6111
$ws.Range("J61").Value = 1299.1666
$ws.Range("K61").Value = 915.6111
$ws.Range("L61").Value = 1299.1666
$ws.Range("M61").Value = -703.6111
$ws.Range("N61").Value = -1723.1666
$ws.Range("H132").Value = 1752.2273
$ws.Range("I132").Value = 1237.45
$ws.Range("J132").Value = 6900
$ws.Range("K132").Value = 3712.35
$ws.Range("L132").Value = 20700
$ws.Range("M132").Value = -1182.35
$ws.Range("N132").Value = -25760
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 1069.0333
$ws.Range("I136").Value = 915.6111
$ws.Range("J136").Value = 1299.1666
$ws.Range("K136").Value = 2746.8333
$ws.Range("L136").Value = 3897.4998
$ws.Range("M136").Value = -196.8332999999998
$ws.Range("N136").Value = -8997.4998
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 49900
$ws.Range("J141").Value = 49900
$ws.Range("L141").Value = 49900
$ws.Range("N141").Value = -60260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2603.8718
$ws.Range("I134").Value = 1305.8334
$ws.Range("J134").Value = 6930.6665
$ws.Range("K134").Value = 3917.5002
$ws.Range("L134").Value = 20791.9995
$ws.Range("M134").Value = -1382.5002
$ws.Range("N134").Value = -25861.9995
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 41331.92
$ws.Range("J138").Value = 41331.92
$ws.Range("L138").Value = 41331.92
$ws.Range("N138").Value = -51611.92
$ws.Range("H140").Value = 49995.715
$ws.Range("J140").Value = 49995.715
$ws.Range("L140").Value = 49995.715
$ws.Range("N140").Value = -60355.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12823242
$ws.Range("I31").Value = 1526.174
$ws.Range("J31").Value = 31254458
$ws.Range("K31").Value = 1526.174
$ws.Range("L31").Value = 31254458
$ws.Range("M31").Value = -1231.174
$ws.Range("N31").Value = -31255048
$ws.Range("H34").Value = 12823242
$ws.Range("I34").Value = 1526.174
$ws.Range("J34").Value = 31254458
$ws.Range("K34").Value = 1526.174
$ws.Range("L34").Value = 31254458
$ws.Range("M34").Value = -1324.174
$ws.Range("N34").Value = -31254862
$ws.Range("H122").Value = 2720.7273
$ws.Range("I122").Value = 1266.6666
$ws.Range("K122").Value = 3799.9998
$ws.Range("M122").Value = -1349.9998
$ws.Range("H134").Value = 3467.58
$ws.Range("I134").Value = 3854.4375
$ws.Range("J134").Value = 2779.8333
$ws.Range("K134").Value = 11563.3125
$ws.Range("L134").Value = 8339.499899999999
$ws.Range("M134").Value = -9028.3125
$ws.Range("N134").Value = -13409.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 147.53334
$ws.Range("J33").Value = 126.44444
$ws.Range("L33").Value = 758.66664
$ws.Range("N33").Value = -1324.66664
$ws.Range("H107").Value = 33163.13
$ws.Range("I107").Value = 466
$ws.Range("J107").Value = 68040.07000000001
$ws.Range("K107").Value = 1398
$ws.Range("L107").Value = 204120.21
$ws.Range("M107").Value = 522
$ws.Range("N107").Value = -207960.21
$ws.Range("H113").Value = 698.0179000000001
$ws.Range("I113").Value = 617.26666
$ws.Range("J113").Value = 1028.3636
$ws.Range("K113").Value = 1851.79998
$ws.Range("L113").Value = 3085.0908
$ws.Range("M113").Value = 318.20002
$ws.Range("N113").Value = -7425.0908
$ws.Range("H132").Value = 1870.3125
$ws.Range("I132").Value = 775
$ws.Range("K132").Value = 6975
$ws.Range("M132").Value = -4445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 17430.428
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 23622.6
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 23622.6
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -27962.6
$ws.Range("H137").Value = 41550
$ws.Range("J137").Value = 41550
$ws.Range("L137").Value = 41550
$ws.Range("N137").Value = -51750

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11116.923
$ws.Range("I132").Value = 13661.1
$ws.Range("J132").Value = 8438.842000000001
$ws.Range("K132").Value = 40983.3
$ws.Range("L132").Value = 25316.526
$ws.Range("M132").Value = -38453.3
$ws.Range("N132").Value = -30376.526
$ws.Range("H136").Value = 2975.75
$ws.Range("I136").Value = 1405.2693
$ws.Range("J136").Value = 7059
$ws.Range("K136").Value = 4215.8079
$ws.Range("L136").Value = 21177
$ws.Range("M136").Value = -1665.8079
$ws.Range("N136").Value = -26277

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2031.902
$ws.Range("I136").Value = 636.1053000000001
$ws.Range("J136").Value = 6111.923
$ws.Range("K136").Value = 1908.3159
$ws.Range("L136").Value = 18335.769
$ws.Range("M136").Value = 641.6840999999999
$ws.Range("N136").Value = -23435.769
